$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.253091096878052
$ws.Range("B1").Value = 1.756467223167419
$ws.Range("C1").Value = 3.163914442062378
$ws.Range("D1").Value = 3.788207054138184
$ws.Range("E1").Value = 1.307068586349487
